# perbaikan pada baris ttl pegawai, dan penambahan tombol unset jabatan aktif
#
# Sets the "Status" (column K) for a set of employees from "Aktif" to
# "Non-Aktif" or "Berhenti", mirroring the underlying data change that backs
# the new "unset active position" button. The sheet layout is:
#   row 1 = title, row 3 = header, rows 4-168 = employee data (one per row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nonAktifRows = @(28, 30, 55, 63, 73, 75, 76, 77, 82, 84, 86, 87, 89, 93, 95, 113, 118, 119, 121, 123, 124, 125, 126, 127, 128, 129, 130, 131)
$berhentiRows = @(103, 104, 106)

foreach ($r in $nonAktifRows) {
    $ws.Range("K$r").Value = "Non-Aktif"
}

foreach ($r in $berhentiRows) {
    $ws.Range("K$r").Value = "Berhenti"
}
